# Apply updated "want to attend" / "lowest price" figures scraped at a later
# run (gh-pages data refresh, commit 456a3b4).
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1105
$ws1.Range("F5").Value  = 176
$ws1.Range("G5").Value  = 139
$ws1.Range("F7").Value  = 217
$ws1.Range("F8").Value  = 392
$ws1.Range("F9").Value  = 1014
$ws1.Range("F11").Value = 508
$ws1.Range("F13").Value = 155
$ws1.Range("F14").Value = 12671
$ws1.Range("F16").Value = 5225
$ws1.Range("F17").Value = 5524

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 63

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1105
$ws4.Range("F6").Value  = 176
$ws4.Range("G6").Value  = 139
$ws4.Range("F8").Value  = 217
$ws4.Range("F9").Value  = 392
$ws4.Range("F10").Value = 1014
$ws4.Range("F12").Value = 508
$ws4.Range("F14").Value = 155
$ws4.Range("F15").Value = 12671
$ws4.Range("F16").Value = 63
$ws4.Range("F19").Value = 5225
$ws4.Range("F20").Value = 5524
